$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-10 (A is a running index, B=tipo, C=pais, D=nombre, E=anio)
$data = @(
    @(0, "lanzadera",    "eeuu",   "saturno v", "1967"),
    @(1, "tripulada",    "rusia",  "soyus",     "1967"),
    @(2, "lanzadera",    "eeuu",   "atlas",     "2002"),
    @(3, "tripulada",    "china",  "shenzou",   "1999"),
    @(4, "no tripulada", "rusia",  "luna I",    "1959"),
    @(5, "no tripulada", "europa", "soho",      "1995"),
    @(6, "lanzadera",    "eeuu",   "zenit II",  "1985"),
    @(7, "tripulada",    "eeuu",   "apolo",     "1966"),
    @(8, "no tripulada", "eeuu",   "mariner x", "1973")
)

# Copy the existing formatting of A2 (bordered/centered/bold style) down to A3:A10
$ws.Range("A2").Copy()
$ws.Range("A3:A10").PasteSpecial(-4122)

# Make sure the E column (years) is treated as text before writing values,
# so strings like "1967" are not auto-converted to numbers.
$ws.Range("E2:E10").NumberFormat = "@"

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}

# Restore the default (Normal) style for the E column so that it doesn't keep
# a lingering text-number-format style, matching the original workbook look.
$ws.Range("E2:E10").Style = "Normal"
